# Update for Vehicle Renewal & Login For Corporate
#
# Sheet "vehicleRenewalTestData" (4th sheet) drops two columns that are no
# longer needed by the test data ("Inspected Status" and
# "Has UAE And GCC AND SALIK Fines"), which also shifts "ExpiredDaysCount"
# and "toRun" left and removes their now-unused shared strings
# (PASSED / FALSE / "Has UAE And GCC AND SALIK Fines").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# Delete column H ("Has UAE And GCC AND SALIK Fines") first, then column F
# ("Inspected Status") - deleting right-to-left keeps the remaining column
# letters stable while each delete is applied.
$ws.Columns("H").Delete()
$ws.Columns("F").Delete()

# Update page setup (paper size / orientation) for the sheet.
$ws.PageSetup.PaperSize = 256
$ws.PageSetup.Orientation = 1

# Move the active selection/cursor as recorded in the saved view state.
[void]$ws.Range("I19").Select()
